$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns for each cryptocurrency row
# D-column values that look like plain decimal numbers must be forced to
# text (NumberFormat "@") before assignment so Excel keeps the exact
# string representation (e.g. preserving trailing zeros) instead of
# silently converting them to floating point numbers.

$ws.Range("D2").Value = "38.307.90"
$ws.Range("E2").Value = "  +3.47%  "

$ws.Range("D3").Value = "2.069.31"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.02"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("E6").Value = "  +2.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.61"
$ws.Range("E7").Value = "  +11.60%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +4.50%  "

$ws.Range("E10").Value = "  +4.76%  "

$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.01"
$ws.Range("E12").Value = "  +6.93%  "

$ws.Range("D13").Value = "2.374.79"
$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.53"
$ws.Range("E14").Value = "  +8.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("E15").Value = "  +4.92%  "

$ws.Range("E16").Value = "  +3.61%  "

$ws.Range("D17").Value = "2.075.40"
$ws.Range("E17").Value = "  +3.17%  "

$ws.Range("D18").Value = "38.234.94"
$ws.Range("E18").Value = "  +3.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.36"
$ws.Range("E20").Value = "  +2.30%  "

$ws.Range("E21").Value = "  +3.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.09"
$ws.Range("E22").Value = "  +1.94%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("E25").Value = "  +3.79%  "

$ws.Range("E26").Value = "  +4.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.73"
$ws.Range("E27").Value = "  +1.34%  "

$ws.Range("E28").Value = "  +5.92%  "

$ws.Range("E29").Value = "  +3.34%  "

$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("E31").Value = "  +3.22%  "

$ws.Range("E32").Value = "  +4.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.67"
$ws.Range("E33").Value = "  +5.21%  "

$ws.Range("E34").Value = "  +9.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0610"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("E37").Value = "  +16.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +6.73%  "

$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").Value = "1.527.52"
$ws.Range("E40").Value = "  +4.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.40"
$ws.Range("E41").Value = "  +9.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.63"
$ws.Range("E42").Value = "  +4.57%  "

$ws.Range("E43").Value = "  +3.88%  "

$ws.Range("E44").Value = "  +4.08%  "

$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.05"
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("E49").Value = "  +2.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.10"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").Value = "2.261.43"
$ws.Range("E51").Value = "  +3.05%  "
